# Insert a new data row at row 733 (pushing existing rows 733-774 down to
# 734-775) and populate it with the new day's data: 2026/02/01 (日), hour 3,
# ranking 22. This mirrors the author's edit which grew the sheet from
# A1:D774 to A1:D775.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 733.. down by one to make room for the new row.
$ws.Rows.Item(733).Insert()

# Column A holds date-like text (e.g. "2026/02/01"). Force it to be stored
# as text (matching the rest of the column) instead of letting Excel
# auto-convert it to a date serial number, then drop the temporary format
# so the cell ends up unstyled like its neighbours.
$ws.Range("A733").NumberFormat = "@"
$ws.Range("A733").Value = "2026/02/01"
$ws.Range("A733").ClearFormats()

$ws.Range("B733").Value = "日"
$ws.Range("C733").Value = 3
$ws.Range("D733").Value = 22
